$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting COR / MARCA PNEU / DIAS ESTOQUE one column right
$ws.Columns("E").Insert()

# New header for inserted column E
$ws.Range("E1").Value = "MODELO"

# Replace the old " 24/24" / "23/24" text values in column D with numeric year 2024,
# and populate the newly inserted column E with the same numeric year value.
$ws.Range("D2").Value = 2024
$ws.Range("E2").Value = 2024

$ws.Range("D3").Value = 2024
$ws.Range("E3").Value = 2024

# Reflect the final saved selection state
$ws.Range("D2").Select()
